$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7000
$ws.Range("J40").Value = 6500
$ws.Range("L40").Value = 6500
$ws.Range("N40").Value = -6850

$ws.Range("H51").Value = 8497.833000000001
$ws.Range("J51").Value = 9999.666999999999
$ws.Range("L51").Value = 9999.666999999999
$ws.Range("N51").Value = -10967.667

$ws.Range("H82").Value = 2049.5
$ws.Range("I82").Value = 2049.5
$ws.Range("K82").Value = 6148.5
$ws.Range("M82").Value = -5742.5

$ws.Range("H85").Value = 2049.5
$ws.Range("I85").Value = 2049.5
$ws.Range("K85").Value = 6148.5
$ws.Range("M85").Value = -4744.5

$ws.Range("H97").Value = 1470
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1470
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 4410
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = -5402

$ws.Range("H99").Value = 325
$ws.Range("I99").Value = 325
$ws.Range("K99").Value = 975
$ws.Range("M99").Value = 523

$ws.Range("H101").Value = 3839
$ws.Range("I101").Value = 3839
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 11517
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -9895
$ws.Range("N101").Value = ""

$ws.Range("H103").Value = 4899.5713
$ws.Range("I103").Value = 4916
$ws.Range("K103").Value = 14748
$ws.Range("M103").Value = -14162

$ws.Range("H107").Value = 622.5
$ws.Range("I107").Value = 603
$ws.Range("K107").Value = 603
$ws.Range("M107").Value = 1317

$ws.Range("H138").Value = 8508
$ws.Range("I138").Value = 7539.6
$ws.Range("K138").Value = 22618.8
$ws.Range("M138").Value = -17478.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8600
$ws.Range("I37").Value = 3000
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 3000
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -2727
$ws.Range("N37").Value = -10546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 18513.834
$ws.Range("J50").Value = 20000
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -21250

$ws.Range("H59").Value = 29999.75
$ws.Range("J59").Value = 29999.75
$ws.Range("L59").Value = 29999.75
$ws.Range("N59").Value = -32289.75

$ws.Range("H68").Value = 85000
$ws.Range("J68").Value = 85000
$ws.Range("L68").Value = 85000
$ws.Range("N68").Value = -86498

$ws.Range("H70").Value = 24999.666
$ws.Range("J70").Value = 24999.666
$ws.Range("L70").Value = 24999.666
$ws.Range("N70").Value = -25629.666

$ws.Range("H71").Value = 85000
$ws.Range("J71").Value = 85000
$ws.Range("L71").Value = 255000
$ws.Range("N71").Value = -262488

$ws.Range("H73").Value = 24999.666
$ws.Range("J73").Value = 24999.666
$ws.Range("L73").Value = 24999.666
$ws.Range("N73").Value = -27183.666

$ws.Range("H132").Value = 1155.5
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1199.3334
$ws.Range("I69").Value = 1199.3334
$ws.Range("K69").Value = 3598.0002
$ws.Range("M69").Value = -2787.0002

$ws.Range("H72").Value = 1199.3334
$ws.Range("I72").Value = 1199.3334
$ws.Range("K72").Value = 10794.0006
$ws.Range("M72").Value = -6738.000599999999

$ws.Range("H86").Value = 1475
$ws.Range("I86").Value = 1475
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4425
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3239
$ws.Range("N86").Value = ""

$ws.Range("H89").Value = 1475
$ws.Range("I89").Value = 1475
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 13275
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -7347
$ws.Range("N89").Value = ""

$ws.Range("H113").Value = 2000
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 6000
$ws.Range("N113").Value = -10340

$ws.Range("H131").Value = 982.0833
$ws.Range("J131").Value = 989.5454999999999
$ws.Range("L131").Value = 2968.6365
$ws.Range("N131").Value = -13048.6365

$ws.Range("H140").Value = 889.5
$ws.Range("I140").Value = 889.5
$ws.Range("K140").Value = 2668.5
$ws.Range("M140").Value = 2511.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12630.857
$ws.Range("I43").Value = 9683.200000000001
$ws.Range("K43").Value = 9683.200000000001
$ws.Range("M43").Value = -9532.200000000001

$ws.Range("H46").Value = 9974.833000000001
$ws.Range("I46").Value = 9924.5
$ws.Range("K46").Value = 9924.5
$ws.Range("M46").Value = -9768.5

$ws.Range("H57").Value = 10680
$ws.Range("I57").Value = 4466.6665
$ws.Range("K57").Value = 4466.6665
$ws.Range("M57").Value = -3646.6665

$ws.Range("H80").Value = 12999.2
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -26996

$ws.Range("H83").Value = 12999.2
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 125000
$ws.Range("N83").Value = -134984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 500
$ws.Range("J21").Value = 500
$ws.Range("L21").Value = 500
$ws.Range("N21").Value = -848

$ws.Range("H55").Value = 1500
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1500
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 1500
$ws.Range("M55").Value = ""
$ws.Range("N55").Value = -1846

$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""

$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 15000
$ws.Range("I21").Value = 20000
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 20000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -19765
$ws.Range("N21").Value = -10470

$ws.Range("H35").Value = 15000
$ws.Range("I35").Value = 20000
$ws.Range("J35").Value = 10000
$ws.Range("K35").Value = 20000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = -19710
$ws.Range("N35").Value = -10580

$ws.Range("H44").Value = 19999
$ws.Range("J44").Value = 19999
$ws.Range("L44").Value = 19999
$ws.Range("N44").Value = -21107

$ws.Range("H69").Value = 10000
$ws.Range("J69").Value = 10000
$ws.Range("L69").Value = 10000
$ws.Range("N69").Value = -11498

$ws.Range("H72").Value = 10000
$ws.Range("J72").Value = 10000
$ws.Range("L72").Value = 30000
$ws.Range("N72").Value = -37488

$ws.Range("H75").Value = 20750
$ws.Range("J75").Value = 20333.334
$ws.Range("L75").Value = 20333.334
$ws.Range("N75").Value = -22205.334

$ws.Range("H78").Value = 20750
$ws.Range("J78").Value = 20333.334
$ws.Range("L78").Value = 61000.00199999999
$ws.Range("N78").Value = -70360.00199999999

$ws.Range("H81").Value = 277
$ws.Range("I81").Value = 277
$ws.Range("K81").Value = 554
$ws.Range("M81").Value = 507

$ws.Range("H84").Value = 277
$ws.Range("I84").Value = 277
$ws.Range("K84").Value = 2770
$ws.Range("M84").Value = 2534

$ws.Range("H130").Value = 49887.332
$ws.Range("J130").Value = 49887.332
$ws.Range("L130").Value = 49887.332
$ws.Range("N130").Value = -59927.332
